# Insert a new weekly price record for "Zapallo italiano" as row 142,
# pushing the existing rows 142-204 down to 143-205 (dimension A1:R204 -> A1:R205).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 142..204 down by one, leaving a blank row 142 ready to fill in.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new record's data.
$ws.Range("A142").Value = 6
$ws.Range("B142").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C142").Value = "Metropolitana"
$ws.Range("D142").Value = 44452
$ws.Range("E142").Value = 13
$ws.Range("F142").Value = 100112032
$ws.Range("G142").Value = "Zapallo italiano"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 330
$ws.Range("K142").Value = 14000
$ws.Range("L142").Value = 15000
$ws.Range("M142").Value = 14697
$ws.Range("N142").Value = "`$/caja 50 unidades"
$ws.Range("O142").Value = "Región de Arica y Parinacota"
$ws.Range("P142").Value = 294
$ws.Range("Q142").Value = 50
$ws.Range("R142").Value = "Hortaliza"
